$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 15
$ws.Range("A6").Value = 4
$ws.Range("A11").Value = 3
$ws.Range("A12").Value = 3
